$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each product row (3..94), "Diferencia Stock" (column L) is corrected
# to match "Stock Minimo Objetivo" (column K).
$total = 0
for ($row = 3; $row -le 94; $row++) {
    $k = $ws.Cells.Item($row, 11).Value()
    $ws.Cells.Item($row, 12).Value = $k
    $total = $total + $k
}

# Update the "Total_Ajuste_Stock" summary cell (C108) to reflect the new sum
# of the "Diferencia Stock" column.
$ws.Range("C108").Value = $total
